$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) figures.
# Values that are plain decimal numbers (a single "." and only digits,
# optionally signed) are written via Formula with a leading apostrophe so
# they stay text cells (matching the original inline-string cells) instead
# of being auto-converted to numeric values by Excel.

$ws.Range("D2").Value = "27.675.48"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "1.584.47"
$ws.Range("E3").Value = "  -3.20%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Formula = "'206.28"
$ws.Range("E5").Value = "  -2.64%  "
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E8").Value = "  -4.86%  "
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("E10").Value = "  -3.15%  "
$ws.Range("E11").Value = "  -1.87%  "
$ws.Range("D12").Value = "1.810.15"
$ws.Range("E12").Value = "  -3.12%  "
$ws.Range("D13").Value = "1.564.11"
$ws.Range("E13").Value = "  -4.46%  "
$ws.Range("E14").Value = "  -4.13%  "
$ws.Range("E15").Value = "  -5.96%  "
$ws.Range("D16").Value = "27.649.23"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").Formula = "'63.29"
$ws.Range("D18").Formula = "'220.13"
$ws.Range("E18").Value = "  -4.01%  "
$ws.Range("E19").Value = "  -3.82%  "
$ws.Range("D20").Formula = "'7.30"
$ws.Range("E20").Value = "  -5.97%  "
$ws.Range("D22").Formula = "'4.13"
$ws.Range("E22").Value = "  -5.46%  "
$ws.Range("E23").Value = "  -6.91%  "
$ws.Range("D24").Formula = "'1.96"
$ws.Range("E24").Value = "  -5.74%  "
$ws.Range("D25").Formula = "'153.85"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("E26").Value = "  -2.83%  "
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("E28").Value = "  -2.84%  "
$ws.Range("E29").Value = "  -4.28%  "
$ws.Range("E30").Value = "  -3.02%  "
$ws.Range("E31").Value = "  -3.56%  "
$ws.Range("E32").Value = "  -6.10%  "
$ws.Range("D33").Value = "1.385.59"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("E34").Value = "  -5.57%  "
$ws.Range("E35").Value = "  -5.41%  "
$ws.Range("D36").Formula = "'0.963"
$ws.Range("E36").Value = "  -5.24%  "
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("E38").Value = "  -3.11%  "
$ws.Range("E39").Value = "  -3.63%  "
$ws.Range("D40").Formula = "'0.819"
$ws.Range("E40").Value = "  -4.08%  "
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").Formula = "'0.978"
$ws.Range("E42").Value = "  -2.99%  "
$ws.Range("D43").Formula = "'1.77"
$ws.Range("E43").Value = "  -4.46%  "
$ws.Range("D44").Formula = "'2.17"
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("D45").Formula = "'63.56"
$ws.Range("E45").Value = "  -3.91%  "
$ws.Range("D46").Formula = "'5.23"
$ws.Range("E46").Value = "  -4.33%  "
$ws.Range("D47").Value = "1.721.43"
$ws.Range("E47").Value = "  -3.10%  "
$ws.Range("D48").Formula = "'87.91"
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").Formula = "'0.0972"
$ws.Range("E50").Value = "  -5.18%  "
$ws.Range("E51").Value = "  -1.04%  "
